# New networks and sampling code
#
# Row 8 previously held the "Cropped_Compensated_Normalized_5FoldPartition"
# label with no metric data yet. This run adds a new "_Test1" sampling
# result: rename the label and fill in its fscore/precision/recall/accuracy
# row, then leave the selection where data entry left off (E9, just below
# the newly-completed row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Cropped_Compensated_Normalized_5FoldPartition_Test1"
$ws.Range("B8").Value = 0.3834
$ws.Range("C8").Value = 0.4607
$ws.Range("D8").Value = 0.5265
$ws.Range("E8").Value = 0.6532

$ws.Range("E9").Select()
